$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks before rewriting data so we can re-add them cleanly
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-12-10 18:28:01'
$ws.Range("B2").Value = '産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5450864'
$ws.Range("G2").Value = 383
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

# Row 3
$ws.Range("A3").Value = '2025-12-10 18:28:01'
$ws.Range("B3").Value = '退職代行サービスの問い合わせや、即時契約直後に聞き取る内容を完了させるAIチャットボットの開発'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '5,000,000 円 ~ / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5451344'
$ws.Range("G3").Value = 375
$ws.Range("H3").Value = '🔥AI,Ai ◆開発'

# Row 4
$ws.Range("A4").Value = '2025-12-10 18:28:01'
$ws.Range("B4").Value = '【自動化】Webサービス更新ツール開発(200アカウント管理)'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5448409'
$ws.Range("G4").Value = 230
$ws.Range("H4").Value = '◆ツール,開発 ◇管理'

# Row 5
$ws.Range("A5").Value = '2025-12-10 18:28:01'
$ws.Range("B5").Value = '【Java/対話システム/心理学実験】協同問題解決プラットフォームの改修開発'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5439921'
$ws.Range("G5").Value = 155
$ws.Range("H5").Value = '★Java ◆開発'

# Row 6
$ws.Range("A6").Value = '2025-12-10 18:28:01'
$ws.Range("B6").Value = 'JavaScriptのスクラッチ開発案件 長期対応可能なパートナー様募集'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5451176'
$ws.Range("G6").Value = 140
$ws.Range("H6").Value = '★Java ◆開発'

# Row 7
$ws.Range("A7").Value = '2025-12-10 18:28:01'
$ws.Range("B7").Value = '海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5251319'
$ws.Range("G7").Value = 135
$ws.Range("H7").Value = '◆ツール,スクレイピング ◇サイト'

# Row 8
$ws.Range("A8").Value = '2025-12-10 18:28:01'
$ws.Range("B8").Value = '【急募】某新聞社のプロトタイプシステム用チャットボット開発'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5450641'
$ws.Range("G8").Value = 83
$ws.Range("H8").Value = '◆開発'

# Row 9
$ws.Range("A9").Value = '2025-12-10 18:28:01'
$ws.Range("B9").Value = 'サーバーサイド保守・追加開発'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5451285'
$ws.Range("G9").Value = 75
$ws.Range("H9").Value = '◆開発'

# Row 10
$ws.Range("A10").Value = '2025-12-10 18:28:01'
$ws.Range("B10").Value = '在宅専業OK│フルスタックエンジニア/開発×データ処理に挑戦!EC運営を支える仕事!'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5450846'
$ws.Range("G10").Value = 75
$ws.Range("H10").Value = '◆開発'

# Row 11
$ws.Range("A11").Value = '2025-12-10 18:28:01'
$ws.Range("B11").Value = '【フルスタックエンジニア募集】新規Webサービス開発'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5450548'
$ws.Range("G11").Value = 75
$ws.Range("H11").Value = '◆開発'

# Row 12
$ws.Range("A12").Value = '2025-12-10 18:28:01'
$ws.Range("B12").Value = '就労継続支援事業所のポータルサイト制作'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5451305'
$ws.Range("G12").Value = 38
$ws.Range("H12").Value = '◇サイト'

# Row 13
$ws.Range("A13").Value = '2025-12-10 18:28:01'
$ws.Range("B13").Value = '【急募】オンラインガチャサイトに決済機能を導入可能な方'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5450884'
$ws.Range("G13").Value = 33
$ws.Range("H13").Value = '◇サイト'

# Row 14
$ws.Range("A14").Value = '2025-12-10 18:28:01'
$ws.Range("B14").Value = '【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5445466'
$ws.Range("G14").Value = 25
$ws.Range("H14").Value = ""

# Row 15
$ws.Range("A15").Value = '2025-12-10 18:28:01'
$ws.Range("B15").Value = '注目 限定公開 PR 限定公開の仕事'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5450323'
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = ""

# Row 16
$ws.Range("A16").Value = '2025-12-10 18:28:01'
$ws.Range("B16").Value = '【急募】当社HPのバグ修正をお願いしたいです'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5450784'
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = ""

# Re-add hyperlinks for F2:F16 in order, matching the URL already placed in each cell
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5450864') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5451344') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5448409') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5439921') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5451176') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5251319') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5450641') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5451285') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5450846') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5450548') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5451305') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5450884') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5445466') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5450323') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5450784') | Out-Null
